$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new line entries (line7, line8) are inserted logically right after
# line6 / before extr1. Every "extr#" row from the old layout shifts down
# by two positions, so every row from 8 through 15 needs both its label
# (column B) and its C/D/E data refreshed, and two new rows (16, 17) are
# appended for extr7 / extr8.

# Row 8: line7
$ws.Range("B8").Value = "line7"
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 11

# Row 9: line8
$ws.Range("B9").Value = "line8"
$ws.Range("C9").Value = 16

# Row 10: extr1
$ws.Range("B10").Value = "extr1"
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 12
$ws.Range("E10").Value = $true

# Row 11: extr2
$ws.Range("B11").Value = "extr2"
$ws.Range("C11").Value = 5
$ws.Range("D11").Value = 9

# Row 12: extr3
$ws.Range("B12").Value = "extr3"
$ws.Range("C12").Value = 10

# Row 13: extr4
$ws.Range("B13").Value = "extr4"
$ws.Range("D13").Value = 8
$ws.Range("E13").Value = $true

# Row 14: extr5
$ws.Range("B14").Value = "extr5"
$ws.Range("C14").Value = 9
$ws.Range("D14").Value = 11

# Row 15: extr6
$ws.Range("B15").Value = "extr6"
$ws.Range("C15").Value = 7
$ws.Range("D15").Value = 11
$ws.Range("E15").Value = $false

# New rows 16 & 17 for extr7 / extr8, carrying forward row 15's
# formatting (bordered, centered, bold) for column A.
$ws.Range("A15").Copy() | Out-Null
$ws.Range("A16:A17").PasteSpecial(-4122) | Out-Null

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "extr7"
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = $true

$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "extr8"
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = $true
